# Update the ISIC code labels in row 10 of the "Pre ISIC Consolidation" sheet.
# These labels act as SUMIFS match criteria for the SoDSCbRIC sheet, so
# changing them re-allocates (disaggregates) the ISIC spending shares.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre ISIC Consolidation")

$ws.Range("E10").Value  = "ISIC 07T08"
$ws.Range("N10").Value  = "ISIC 22"
$ws.Range("S10").Value  = "ISIC 25"
$ws.Range("T10").Value  = "ISIC 26"
$ws.Range("U10").Value  = "ISIC 27"
$ws.Range("X10").Value  = "ISIC 30"

# Recalculate the workbook so dependent formulas (e.g. the SUMIFS on the
# SoDSCbRIC sheet) refresh their cached results.
$excel.CalculateFullRebuild()

# Leave the "Pre ISIC Consolidation" sheet active/selected, matching the
# author's final view state.
$ws.Activate() | Out-Null
$ws.Range("AA23").Select() | Out-Null
